$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,7).Value = 1.38698
$ws.Cells.Item(2,8).Value = 2.77396
$ws.Cells.Item(2,9).Value = 0.5967131687073423
$ws.Cells.Item(2,10).Value = 0.4965808999056411
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,13).Value = 0.6483006666666666
$ws.Cells.Item(2,14).Value = 1.944902
$ws.Cells.Item(2,15).Value = 0.2310664473662325
$ws.Cells.Item(2,16).Value = 0.3059082590740972
$ws.Cells.Item(2,17).Value = 0.8991800586533332
$ws.Cells.Item(2,18).Value = 5.395080351919999
$ws.Cells.Item(2,19).Value = 0.1378803919898529
$ws.Cells.Item(2,20).Value = 0.1519081985795832
$ws.Cells.Item(3,4).Value = "M2"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,7).Value = 1.38698
$ws.Cells.Item(3,8).Value = 2.77396
$ws.Cells.Item(3,9).Value = 0.5967131687073423
$ws.Cells.Item(3,10).Value = 0.4965808999056411
$ws.Cells.Item(3,12).Value = 0.3333333333333333
$ws.Cells.Item(3,13).Value = 0.029863
$ws.Cells.Item(3,14).Value = 0.089589
$ws.Cells.Item(3,15).Value = 0.01064373009698864
$ws.Cells.Item(3,16).Value = 0.01409120614930176
$ws.Cells.Item(3,17).Value = 0.04141938374
$ws.Cells.Item(3,18).Value = 0.24851630244
$ws.Cells.Item(3,19).Value = 0.006351253913039799
$ws.Cells.Item(3,20).Value = 0.006997423830376172
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,4).Value = "Neutro"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,7).Value = 1.38698
$ws.Cells.Item(4,8).Value = 2.77396
$ws.Cells.Item(4,9).Value = 0.5967131687073423
$ws.Cells.Item(4,10).Value = 0.4965808999056411
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.06825266666666667
$ws.Cells.Item(4,14).Value = 0.204758
$ws.Cells.Item(4,15).Value = 0.0243265232026164
$ws.Cells.Item(4,16).Value = 0.03220581978500407
$ws.Cells.Item(4,17).Value = 0.09466508361333333
$ws.Cells.Item(4,18).Value = 0.56799050168
$ws.Cells.Item(4,19).Value = 0.01451595674386592
$ws.Cells.Item(4,20).Value = 0.01599279497103622
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,7).Value = 1.38698
$ws.Cells.Item(5,8).Value = 2.77396
$ws.Cells.Item(5,9).Value = 0.5967131687073423
$ws.Cells.Item(5,10).Value = 0.4965808999056411
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,13).Value = 2.059273
$ws.Cells.Item(5,14).Value = 4.118546
$ws.Cells.Item(5,15).Value = 0.7339632993341625
$ws.Cells.Item(5,16).Value = 0.6477947149915969
$ws.Cells.Item(5,17).Value = 2.85617046554
$ws.Cells.Item(5,18).Value = 11.42468186216
$ws.Cells.Item(5,19).Value = 0.4379655660605836
$ws.Cells.Item(5,20).Value = 0.3216824825246455
$ws.Cells.Item(6,1).Value = "Neutro"
$ws.Cells.Item(6,2).Value = "Vip"
$ws.Cells.Item(6,3).Value = "Vipr2"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.9373863333333334
$ws.Cells.Item(6,8).Value = 2.812159
$ws.Cells.Item(6,9).Value = 0.4032868312926577
$ws.Cells.Item(6,10).Value = 0.503419100094359
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.6483006666666666
$ws.Cells.Item(6,14).Value = 1.944902
$ws.Cells.Item(6,15).Value = 0.2310664473662325
$ws.Cells.Item(6,16).Value = 0.3059082590740972
$ws.Cells.Item(6,17).Value = 0.6077081848242223
$ws.Cells.Item(6,18).Value = 5.469373663418001
$ws.Cells.Item(6,19).Value = 0.09318605537637958
$ws.Cells.Item(6,20).Value = 0.154000060494514
$ws.Cells.Item(7,1).Value = "Neutro"
$ws.Cells.Item(7,2).Value = "Vip"
$ws.Cells.Item(7,3).Value = "Vipr2"
$ws.Cells.Item(7,4).Value = "M2"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.9373863333333334
$ws.Cells.Item(7,8).Value = 2.812159
$ws.Cells.Item(7,9).Value = 0.4032868312926577
$ws.Cells.Item(7,10).Value = 0.503419100094359
$ws.Cells.Item(7,11).Value = 1
$ws.Cells.Item(7,12).Value = 0.3333333333333333
$ws.Cells.Item(7,13).Value = 0.029863
$ws.Cells.Item(7,14).Value = 0.089589
$ws.Cells.Item(7,15).Value = 0.01064373009698864
$ws.Cells.Item(7,16).Value = 0.01409120614930176
$ws.Cells.Item(7,17).Value = 0.02799316807233334
$ws.Cells.Item(7,18).Value = 0.251938512651
$ws.Cells.Item(7,19).Value = 0.004292476183948842
$ws.Cells.Item(7,20).Value = 0.007093782318925591
$ws.Cells.Item(8,1).Value = "Neutro"
$ws.Cells.Item(8,2).Value = "Vip"
$ws.Cells.Item(8,3).Value = "Vipr2"
$ws.Cells.Item(8,4).Value = "Neutro"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.9373863333333334
$ws.Cells.Item(8,8).Value = 2.812159
$ws.Cells.Item(8,9).Value = 0.4032868312926577
$ws.Cells.Item(8,10).Value = 0.503419100094359
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.06825266666666667
$ws.Cells.Item(8,14).Value = 0.204758
$ws.Cells.Item(8,15).Value = 0.0243265232026164
$ws.Cells.Item(8,16).Value = 0.03220581978500407
$ws.Cells.Item(8,17).Value = 0.0639791169468889
$ws.Cells.Item(8,18).Value = 0.5758120525220001
$ws.Cells.Item(8,19).Value = 0.009810566458750484
$ws.Cells.Item(8,20).Value = 0.01621302481396785
$ws.Cells.Item(9,1).Value = "Neutro"
$ws.Cells.Item(9,2).Value = "Vip"
$ws.Cells.Item(9,3).Value = "Vipr2"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.9373863333333334
$ws.Cells.Item(9,8).Value = 2.812159
$ws.Cells.Item(9,9).Value = 0.4032868312926577
$ws.Cells.Item(9,10).Value = 0.503419100094359
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 2.059273
$ws.Cells.Item(9,14).Value = 4.118546
$ws.Cells.Item(9,15).Value = 0.7339632993341625
$ws.Cells.Item(9,16).Value = 0.6477947149915969
$ws.Cells.Item(9,17).Value = 1.930334366802334
$ws.Cells.Item(9,18).Value = 11.582006200814
$ws.Cells.Item(9,19).Value = 0.2959977332735789
$ws.Cells.Item(9,20).Value = 0.3261122324669515
